$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ">  H1-H6" / MountingHole row (old row 3) -- shifts rows 4-14 up to 3-13
$ws.Rows(3).Delete()

# Row 2: Capacitors C1, C2 -- update Value / Manufacturer Part number
$ws.Range("B2").Value = "Capacitor 10uF"
$ws.Range("C2").Value = "CL21A106KPFNNNF"

# Row 11 (was row 12 "R1"): update Value, clear stray part-number column
$ws.Range("B11").Value = "Resistor 0603 560Ohm"
$ws.Range("C11").ClearContents()

Write-Host "Done"
